$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-CellText($rowIndex, $oldText, $newText) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    # wdFindStop (0) + wdReplaceOne (1) so the match/replace stays scoped
    # to this cell's Range instead of spilling into the rest of the table.
    $cell.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 1)
}

# Simple single-value rows (Times New Roman / sz 22 runs, one <w:t> each).
Replace-CellText 1 "100" "0M"
Replace-CellText 2 "0.01" "0M"
Replace-CellText 3 "2445" "0M"
Replace-CellText 4 "3" "149"
Replace-CellText 5 "0.00005" "0.00003"
Replace-CellText 7 "0.00006" "0.00005"
Replace-CellText 8 "0.00002" "0.00001"
Replace-CellText 12 "0.00019" "0.00733"

# Rows 44-46 currently hold a whole tab-separated record crammed into a
# single cell (several <w:t>/<w:tab/> runs). Collapse each down to a
# single plain value, reusing the cell's existing run formatting.
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.01"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "2445"
